$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-10"

# Update the header label for the current-year column (I1)
$ws.Range("I1").Value = "2022 (through 04-10)"

# Update the April 2022 carjacking count (I5)
$ws.Range("I5").Value = 38

# Update the Total 2022 carjacking count (I14)
$ws.Range("I14").Value = 472
